# Refactor peak detection algorithm for improved accuracy
#
# The updated peak-detection pass recomputes segment durations per
# recording. Row 3 (14_05_Benjamin) gets 4 newly-detected peaks inserted
# in the middle of its segment list (shifting later durations to the
# right) plus 2 new trailing durations; row 2 (05_06_Roxanne) just gets
# its trailing two durations corrected and 2 new ones appended. Both rows
# therefore grow from 73 to 77 segment columns, so the header row (row 1)
# gains 4 new "Segment 74".."Segment 77" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing cells whose fill/border formatting represent the three
# alternating "segment duration" cell styles used throughout the sheet;
# reuse them as live templates instead of hard-coding colors.
$styleOrange = $ws.Range("B2")
$styleGrey   = $ws.Range("C2")
$styleGreen  = $ws.Range("T2")

function Copy-SegmentStyle($targetAddress, $styleCode) {
    $template = $styleGrey
    if ($styleCode -eq "1") { $template = $styleOrange }
    elseif ($styleCode -eq "2") { $template = $styleGrey }
    elseif ($styleCode -eq "3") { $template = $styleGreen }
    $target = $ws.Range($targetAddress)
    $target.Interior.Color = $template.Interior.Color
    $target.Borders.LineStyle = $template.Borders.LineStyle
    $target.Borders.Color = $template.Borders.Color
}

# --- Row 1: four new segment headers appended after BV1 ---
$ws.Range("BW1").Value = "Segment 74"
$ws.Range("BX1").Value = "Segment 75"
$ws.Range("BY1").Value = "Segment 76"
$ws.Range("BZ1").Value = "Segment 77"

# --- Row 2: recomputed trailing durations (BU2:BX2) ---
$row2Values = New-Object "object[,]" 1,4
$row2Values[0,0] = 2.6
$row2Values[0,1] = 0.8
$row2Values[0,2] = 1.6
$row2Values[0,3] = 0.4
$ws.Range("BU2:BX2").Value = $row2Values

$row2Styles = @(
"2", "3", "2", "3"
)
$row2Cells = @(
"BU2", "BV2", "BW2", "BX2"
)
for ($i = 0; $i -lt $row2Cells.Length; $i++) {
    Copy-SegmentStyle $row2Cells[$i] $row2Styles[$i]
}

# --- Row 3: full set of recomputed durations (B3:BZ3) ---
$row3Values = New-Object "object[,]" 1,77
$row3Values[0,0] = 1
$row3Values[0,1] = 7.7
$row3Values[0,2] = 2.8
$row3Values[0,3] = 5
$row3Values[0,4] = 8.2
$row3Values[0,5] = 10.5
$row3Values[0,6] = 2.8
$row3Values[0,7] = 8.1
$row3Values[0,8] = 1
$row3Values[0,9] = 1.4
$row3Values[0,10] = 1.8
$row3Values[0,11] = 7.9
$row3Values[0,12] = 3.8
$row3Values[0,13] = 10.5
$row3Values[0,14] = 4.4
$row3Values[0,15] = 15.1
$row3Values[0,16] = 2.2
$row3Values[0,17] = 11.7
$row3Values[0,18] = 0.8
$row3Values[0,19] = 1
$row3Values[0,20] = 3.6
$row3Values[0,21] = 9.1
$row3Values[0,22] = 4.6
$row3Values[0,23] = 16.2
$row3Values[0,24] = 1.4
$row3Values[0,25] = 9.1
$row3Values[0,26] = 1.4
$row3Values[0,27] = 12.1
$row3Values[0,28] = 2.6
$row3Values[0,29] = 5.2
$row3Values[0,30] = 6.7
$row3Values[0,31] = 3.2
$row3Values[0,32] = 5
$row3Values[0,33] = 12.3
$row3Values[0,34] = 4.8
$row3Values[0,35] = 7.7
$row3Values[0,36] = 2.2
$row3Values[0,37] = 9.5
$row3Values[0,38] = 3.2
$row3Values[0,39] = 1.8
$row3Values[0,40] = 0.2
$row3Values[0,41] = 2.4
$row3Values[0,42] = 4.4
$row3Values[0,43] = 13.1
$row3Values[0,44] = 3
$row3Values[0,45] = 11.1
$row3Values[0,46] = 4
$row3Values[0,47] = 18.6
$row3Values[0,48] = 3.4
$row3Values[0,49] = 14.7
$row3Values[0,50] = 9.5
$row3Values[0,51] = 8.9
$row3Values[0,52] = 3.4
$row3Values[0,53] = 5.5
$row3Values[0,54] = 4.8
$row3Values[0,55] = 8.3
$row3Values[0,56] = 2
$row3Values[0,57] = 10.1
$row3Values[0,58] = 3.4
$row3Values[0,59] = 3.6
$row3Values[0,60] = 7.5
$row3Values[0,61] = 6.3
$row3Values[0,62] = 8.1
$row3Values[0,63] = 15.3
$row3Values[0,64] = 3.8
$row3Values[0,65] = 3.4
$row3Values[0,66] = 5.7
$row3Values[0,67] = 2.6
$row3Values[0,68] = 11.9
$row3Values[0,69] = 1.8
$row3Values[0,70] = 6.7
$row3Values[0,71] = 2.2
$row3Values[0,72] = 3.4
$row3Values[0,73] = 3
$row3Values[0,74] = 4.6
$row3Values[0,75] = 3.1
$row3Values[0,76] = 4.6
$ws.Range("B3:BZ3").Value = $row3Values

$row3Styles = @(
"3", "2", "1", "2", "1", "2", "1", "2", "1", "2", "3", "2", "1", "2", "1", "2", "1", "2", "1", "2", "3", "2", "1", "2", "1", "2", "1", "2", "1", "2", "1", "2", "1", "2", "1", "2", "1", "2", "1", "2", "3", "2", "1", "2", "1", "2", "1", "2", "1", "2", "1", "2", "1", "2", "1", "2", "1", "2", "1", "2", "1", "2", "1", "2", "3", "2", "1", "2", "1", "2", "1", "2", "3", "2", "3", "2", "3"
)
$row3Cells = @(
"B3", "C3", "D3", "E3", "F3", "G3", "H3", "I3", "J3", "K3", "L3", "M3", "N3", "O3", "P3", "Q3", "R3", "S3", "T3", "U3", "V3", "W3", "X3", "Y3", "Z3", "AA3", "AB3", "AC3", "AD3", "AE3", "AF3", "AG3", "AH3", "AI3", "AJ3", "AK3", "AL3", "AM3", "AN3", "AO3", "AP3", "AQ3", "AR3", "AS3", "AT3", "AU3", "AV3", "AW3", "AX3", "AY3", "AZ3", "BA3", "BB3", "BC3", "BD3", "BE3", "BF3", "BG3", "BH3", "BI3", "BJ3", "BK3", "BL3", "BM3", "BN3", "BO3", "BP3", "BQ3", "BR3", "BS3", "BT3", "BU3", "BV3", "BW3", "BX3", "BY3", "BZ3"
)
for ($i = 0; $i -lt $row3Cells.Length; $i++) {
    Copy-SegmentStyle $row3Cells[$i] $row3Styles[$i]
}

# --- Restore the author's last selection ---
$ws.Range("H17").Select()

